$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.133.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.631.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'216.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.97%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.16%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.31%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'20.10"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.70%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.859.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.629.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.29%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.47%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.36%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.09%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'27.111.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.35%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'  -1.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'214.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.96%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.00%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.98%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -1.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'147.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.98%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E27").Value = "'  -0.59%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.10%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.34%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.32%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.73%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.17%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -1.07%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.307.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.47%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -0.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -1.64%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.542"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.844"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.25%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.02%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.22%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.804"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.90%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.03%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.769.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'62.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'90.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.80%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.08%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.813"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +20.40%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0514"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.27%  "
$ws.Range("E51").Style = "Normal"
